# Updated cryptos list (Price + Volume(1h) columns) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that would otherwise be auto-parsed as numbers by Excel need to be
# forced to Text format first so the values round-trip as strings, matching
# the inline-string cells already present in the workbook.
$numericRiskCells = @(
    "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D15",
    "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27",
    "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37",
    "D38", "D40", "D43", "D44", "D47", "D48", "D49", "D50", "D51"
)
foreach ($ref in $numericRiskCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated Price (D) / Volume(1h) (E) values
$ws.Range("D2").Value = '27.312.14'
$ws.Range("E2").Value = '  -1.43%  '
$ws.Range("D3").Value = '1.707.93'
$ws.Range("E3").Value = '  -1.52%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '224.17'
$ws.Range("E5").Value = '  -1.54%  '
$ws.Range("D6").Value = '0.5333'
$ws.Range("E6").Value = '  -2.29%  '
$ws.Range("D7").Value = '1.003'
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '0.2676'
$ws.Range("E8").Value = '  -2.47%  '
$ws.Range("D9").Value = '0.06610'
$ws.Range("E9").Value = '  -1.60%  '
$ws.Range("D10").Value = '20.96'
$ws.Range("E10").Value = '  -4.26%  '
$ws.Range("D11").Value = '0.07625'
$ws.Range("E11").Value = '  -1.94%  '
$ws.Range("D12").Value = '4.552'
$ws.Range("E12").Value = '  -3.12%  '
$ws.Range("D13").Value = '1.721.06'
$ws.Range("E13").Value = '  -0.62%  '
$ws.Range("D14").Value = '1.945.78'
$ws.Range("E14").Value = '  -1.37%  '
$ws.Range("D15").Value = '0.5772'
$ws.Range("E15").Value = '  -3.64%  '
$ws.Range("D16").Value = '0.0₅8179'
$ws.Range("E16").Value = '  -2.82%  '
$ws.Range("D17").Value = '67.76'
$ws.Range("E17").Value = '  -2.05%  '
$ws.Range("D18").Value = '27.344.68'
$ws.Range("E18").Value = '  -1.33%  '
$ws.Range("D19").Value = '217.61'
$ws.Range("E19").Value = '  -4.14%  '
$ws.Range("D20").Value = '1.003'
$ws.Range("E20").Value = '  +0.04%  '
$ws.Range("D21").Value = '4.669'
$ws.Range("E21").Value = '  -3.46%  '
$ws.Range("D22").Value = '10.45'
$ws.Range("E22").Value = '  -4.11%  '
$ws.Range("D23").Value = '5.968'
$ws.Range("E23").Value = '  -4.09%  '
$ws.Range("D24").Value = '1.004'
$ws.Range("D25").Value = '142.07'
$ws.Range("E25").Value = '  -4.11%  '
$ws.Range("D26").Value = '1.739'
$ws.Range("E26").Value = '  +0.82%  '
$ws.Range("D27").Value = '0.1213'
$ws.Range("E27").Value = '  -3.03%  '
$ws.Range("D28").Value = '7.264'
$ws.Range("E28").Value = '  -2.80%  '
$ws.Range("D29").Value = '16.27'
$ws.Range("E29").Value = '  -4.96%  '
$ws.Range("D30").Value = '0.05399'
$ws.Range("E30").Value = '  -5.50%  '
$ws.Range("D31").Value = '1.291'
$ws.Range("E31").Value = '  -1.67%  '
$ws.Range("D32").Value = '3.499'
$ws.Range("E32").Value = '  -5.59%  '
$ws.Range("D33").Value = '3.426'
$ws.Range("E33").Value = '  -2.68%  '
$ws.Range("D34").Value = '1.644'
$ws.Range("E34").Value = '  -2.52%  '
$ws.Range("D35").Value = '2.875'
$ws.Range("E35").Value = '  +0.72%  '
$ws.Range("D36").Value = '0.9490'
$ws.Range("E36").Value = '  -2.82%  '
$ws.Range("D37").Value = '2.414'
$ws.Range("E37").Value = '  -0.95%  '
$ws.Range("D38").Value = '0.5864'
$ws.Range("E38").Value = '  -2.21%  '
$ws.Range("E39").Value = '  -2.02%  '
$ws.Range("D40").Value = '5.853'
$ws.Range("E40").Value = '  -1.16%  '
$ws.Range("D41").Value = '1.045.98'
$ws.Range("E41").Value = '  -0.25%  '
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("D43").Value = '0.8404'
$ws.Range("E43").Value = '  -1.23%  '
$ws.Range("D44").Value = '101.02'
$ws.Range("E44").Value = '  -0.66%  '
$ws.Range("D45").Value = '1.853.30'
$ws.Range("E45").Value = '  -1.30%  '
$ws.Range("E46").Value = '  +2.09%  '
$ws.Range("D47").Value = '58.01'
$ws.Range("E47").Value = '  -2.89%  '
$ws.Range("D48").Value = '0.4510'
$ws.Range("E48").Value = '  +1.81%  '
$ws.Range("D49").Value = '1.008'
$ws.Range("E49").Value = '  +0.16%  '
$ws.Range("D50").Value = '8.082'
$ws.Range("E50").Value = '  -2.29%  '
$ws.Range("D51").Value = '0.05233'
$ws.Range("E51").Value = '  -1.93%  '

